$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.642214298248291
$ws.Range("B1").Value = 4.576918125152588
$ws.Range("C1").Value = 3.535986185073853
$ws.Range("D1").Value = 3.363679170608521
$ws.Range("E1").Value = 2.124365329742432
